$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("apiTest")
$ws.Range("A3").Value = "srdjan.rados@htecgroup"
$null = $ws.Range("A4").Select()
